$d = $word.ActiveDocument

function Insert-ParaXml($range, $innerXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# Find the two placeholder "-" paragraphs that sit right after the
# "Functional Requirements" / "Non-Functional Requirements" headings.
# ---------------------------------------------------------------------------
$funcHeadingRange = $d.Content
$funcHeadingRange.Find.Execute("Functional Requirements", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$funcParaIndex = -1
$nonFuncParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -eq "Non-Functional Requirements`r") {
        $nonFuncParaIndex = $i
    } elseif ($txt -eq "Functional Requirements`r") {
        $funcParaIndex = $i
    }
}

Write-Host "Functional Requirements heading paragraph: " $funcParaIndex
Write-Host "Non-Functional Requirements heading paragraph: " $nonFuncParaIndex

# The placeholder paragraphs are the very next paragraph after each heading.
$funcDashIndex = $funcParaIndex + 1
$nonFuncDashIndex = $nonFuncParaIndex + 1

# ---------------------------------------------------------------------------
# Process the LATER block first (Non-Functional / "shall") so that the
# paragraph indices for the earlier block stay valid while we work.
# ---------------------------------------------------------------------------
$shallPara = $d.Paragraphs.Item($nonFuncDashIndex)
$shallRange = $shallPara.Range
$shallInner = @'
<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The suggested Sensor Fusion Framework </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>shall</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Follow the Objected Oriented Programming (OOP) Design Pattern.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Be easy to extend or modify.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Support the ROS2 (C++).</w:t></w:r></w:p>
'@
Insert-ParaXml $shallRange $shallInner
Write-Host "Inserted shall block"

# ---------------------------------------------------------------------------
# Now process the EARLIER block (Functional / "must").
# ---------------------------------------------------------------------------
$mustPara = $d.Paragraphs.Item($funcDashIndex)
$mustRange = $mustPara.Range
$mustInner = @'
<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The suggested Sensor Fusion Framework </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>must</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> be able to:</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>A</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ccept different sensor architectures.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Allow different fusion methods </w:t></w:r></w:p>
'@
Insert-ParaXml $mustRange $mustInner
Write-Host "Inserted must block"

# ---------------------------------------------------------------------------
# Apply numbering (numId 2 for the "must" bullets, numId 3 for the "shall"
# bullets) by cloning the existing numbered-list template already used by
# the Milestones section (numId 1). Word/IronDocx will mint fresh num
# entries (2, then 3) in numbering.xml automatically.
# ---------------------------------------------------------------------------
$templateSource = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.ListFormat.ListType -ne 0) {
        $templateSource = $pp.Range.ListFormat.ListTemplate
        break
    }
}

# Re-locate the "must" bullet paragraphs by text and number them (numId 2).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -eq "Accept different sensor architectures.`r" -or $txt -eq "Allow different fusion methods `r") {
        $d.Paragraphs.Item($i).Range.ListFormat.ApplyListTemplateWithLevel($templateSource)
    }
}

# Re-locate the "shall" bullet paragraphs and number them (numId 3).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -eq "Follow the Objected Oriented Programming (OOP) Design Pattern.`r" -or $txt -eq "Be easy to extend or modify.`r" -or $txt -eq "Support the ROS2 (C++).`r") {
        $d.Paragraphs.Item($i).Range.ListFormat.ApplyListTemplateWithLevel($templateSource)
    }
}

Write-Host "Done applying numbering"
